$wb = $excel.ActiveWorkbook

$wsEntries = $wb.Worksheets.Item("Entries")
$wsRules = $wb.Worksheets.Item("Rules")

# Insert a new "Entries" row before the current row 7 (Deny all access),
# shifting the existing rows 7-9 down to 8-10, and populate it with the
# new "send audit record" entry.
$wsEntries.Rows.Item(7).Insert()
$wsEntries.Range("A7").Value = "Audit denied other access and notify and send"
$wsEntries.Range("B7").Value = "AuditDenied"
$wsEntries.Range("D7").Value = "X"
$wsEntries.Range("E7").Value = "X"
$wsEntries.Range("J7").Value = 3

# Update the "Rules" sheet text to reflect the new wording for the
# "send audit record" variant.
$wsRules.Range("E3").Value = "Deny other access, Audit denied other access and notify and send"

# Update view/selection state: Entries keeps its old selection (it is no
# longer the active tab), while Rules becomes the active tab/sheet with a
# new selection.
$wsEntries.Range("A7").Select()
$wsRules.Activate()
$wsRules.Range("E4").Select()
